$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-12-30 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-12-31 Wednesday", 2)
$d.Content.Find.Execute("575÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "110÷6=", 2)
$d.Content.Find.Execute("174÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "863÷8=", 2)
$d.Content.Find.Execute("841÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "919÷5=", 2)
$d.Content.Find.Execute("513÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "181÷4=", 2)
$d.Content.Find.Execute("310÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "829÷9=", 2)
$d.Content.Find.Execute("632÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "898÷5=", 2)
$d.Content.Find.Execute("640÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "696÷8=", 2)
$d.Content.Find.Execute("151÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "364÷9=", 2)
$d.Content.Find.Execute("100÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "965÷7=", 2)
$d.Content.Find.Execute("983÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "901÷5=", 2)
$d.Content.Find.Execute("270÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "702÷9=", 2)
$d.Content.Find.Execute("731÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "997÷2=", 2)
$d.Content.Find.Execute("731÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "364÷2=", 2)
$d.Content.Find.Execute("247÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "898÷7=", 2)
$d.Content.Find.Execute("779÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "364÷3=", 2)
$d.Content.Find.Execute("454÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "627÷6=", 2)
$d.Content.Find.Execute("629÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "618÷7=", 2)
$d.Content.Find.Execute("854÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "525÷2=", 2)
$d.Content.Find.Execute("747÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "978÷3=", 2)
$d.Content.Find.Execute("636÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "623÷8=", 2)
$d.Content.Find.Execute("980÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "451÷9=", 2)
$d.Content.Find.Execute("600÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "248÷5=", 2)
$d.Content.Find.Execute("251÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "103÷9=", 2)
$d.Content.Find.Execute("900÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "153÷7=", 2)
$d.Content.Find.Execute("630÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "805÷3=", 2)
